$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MonAn")

# New dish record (id 12) added by "coder thanh son"
$name = "Dưa góp đu đủ cà rốt"
$ingredients = "1 quả đu đủ" + [char]10 + "1 củ cà rốt" + [char]10 + "4 tép tỏi" + [char]10 + "1 quả ớt băm nhỏ" + [char]10 + "Nước mắm" + [char]10 + "Gia vị khác"
$prep = "Đu đủ nạo vỏ, hạt trong ruột, ngâm nước cho ra hết nhựa.Cà rốt rửa sạch, thái lát mỏng, vừa ăn"
$usage = "Món dưa góp muối nước mắm này để khoảng 1 ngày là ăn được, có thể để cỡ 1 tuần nếu cho vào ngăn mát tủ lạnh."
$cook = "Cho 2 thìa nhỏ muối và 2 thìa canh đường vào bát sau đó trộn đều để ngấm.Đu đủ, cà rốt, dưa chuột trộn đều rồi cho vào hũ thủy tinh đã tiệt trùng  rồi chế nước ngâm vào sao cho ngập hỗn hợp nguyên liệu. Nhớ gài que tre hay bát đĩa để đu đủ, cà rốt, dưa chuột không nổi lên rồi đậy kín hũ. Để hủ nơi thoáng mát."

$row = 13

$ws.Cells.Item($row, 1).Value = 12
$ws.Cells.Item($row, 2).Value = $name
$ws.Cells.Item($row, 6).Value = $ingredients
$ws.Cells.Item($row, 7).Value = $prep
$ws.Cells.Item($row, 8).Value = $cook
$ws.Cells.Item($row, 9).Value = $usage
$ws.Cells.Item($row, 10).Value = 1

# Give the dish name its own bold, small, dark-grey Segoe UI look
$nameCell = $ws.Cells.Item($row, 2)
$nameCell.ClearFormats()
$nameCell.Font.Bold = $true
$nameCell.Font.Size = 9
$nameCell.Font.Name = "Segoe UI"
$nameCell.Font.Color = 3815994

$ws.Rows.Item($row).RowHeight = 129.6

$ws.PageSetup.Orientation = 1

$ws.Range("F20").Select()
